$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple text/value updates (strings that Excel will not mis-parse as numbers)
$ws.Range("D2").Value = '27.163.44'
$ws.Range("E2").Value = '  -2.16%  '
$ws.Range("D3").Value = '1.568.94'
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("E6").Value = '  -2.81%  '
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -2.29%  '
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D12").Value = '1.790.25'
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("D13").Value = '1.561.90'
$ws.Range("E13").Value = '  -2.79%  '
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("E15").Value = '  -2.53%  '
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").Value = '27.161.12'
$ws.Range("E17").Value = '  -2.18%  '
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("E19").Value = '  -1.25%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0686'
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("E23").Value = '  -3.60%  '
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("E26").Value = '  -7.34%  '
$ws.Range("E27").Value = '  -1.49%  '
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("E29").Value = '  -1.54%  '
$ws.Range("E30").Value = '  -2.92%  '
$ws.Range("E31").Value = '  -2.29%  '
$ws.Range("E32").Value = '  -1.95%  '
$ws.Range("D33").Value = '1.396.92'
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("E34").Value = '  -1.73%  '
$ws.Range("E35").Value = '  +0.95%  '
$ws.Range("E37").Value = '  -3.71%  '
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("E39").Value = '  -1.64%  '
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("E42").Value = '  +1.72%  '
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("E44").Value = '  +2.01%  '
$ws.Range("E45").Value = '  -1.38%  '
$ws.Range("E46").Value = '  -0.27%  '
$ws.Range("D47").Value = '1.702.66'
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("D49").Value = '0.0₇0979'
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E50").Value = '  -1.81%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("E51").Value = '  -0.56%  '

# Numeric-looking Price strings must stay TEXT (source data is plain text),
# so force text formatting, assign, then restore default style to avoid
# leaving a visible number-format override on the cell.
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '206.90'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.488'
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '22.35'
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0591'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.520'
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '62.97'
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '214.83'
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '7.30'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.39'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '152.30'
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '6.66'
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '14.96'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.13'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.0463'
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.18'
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.56'
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0166'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.816'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.80'
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.34'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '63.62'
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '85.83'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0950'
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0493'
$c.Style = "Normal"
